# Apply cryptos list update (GitHub Actions scheduled refresh, Thu Nov 16 08:44:57 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.392.57"
$ws.Range("E2").Value = "  +4.41%  "

# Row 3
$ws.Range("D3").Value = "2.046.09"
$ws.Range("E3").Value = "  +2.86%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "253.27"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.10%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.653"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.86%  "

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "65.78"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +10.61%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.409"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +11.99%  "

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "59.81"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.63%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0782"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +5.29%  "

# Row 12
$ws.Range("E12").Value = "  +0.05%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.928"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.64%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "23.66"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +26.01%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "14.84"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.54%  "

# Row 16
$ws.Range("D16").Value = "2.346.92"
$ws.Range("E16").Value = "  +2.99%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.74"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +7.49%  "

# Row 18
$ws.Range("D18").Value = "2.046.89"
$ws.Range("E18").Value = "  +3.05%  "

# Row 19
$ws.Range("D19").Value = "37.261.72"
$ws.Range("E19").Value = "  +4.25%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "73.63"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +2.43%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  +4.02%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.53"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +5.97%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "240.44"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.93%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.63"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "

# Row 26
$ws.Range("E26").Value = "  +4.72%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.18"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +9.20%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "161.95"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.74%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "20.12"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.91%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.133"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +35.65%  "

# Row 31
$ws.Range("E31").Value = "  +2.86%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.20"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +5.06%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.20"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +5.50%  "

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0632"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +4.80%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.71"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +6.96%  "

# Row 36
$ws.Range("E36").Value = "  -3.06%  "

# Row 37
$ws.Range("E37").Value = "  +10.89%  "

# Row 38
$ws.Range("E38").Value = "  +0.06%  "

# Row 39
$ws.Range("E39").Value = "  +2.84%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.07"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +31.90%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +8.34%  "

# Row 42
$ws.Range("E42").Value = "  +3.37%  "

# Row 43
$ws.Range("E43").Value = "  +7.19%  "

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "17.83"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +7.94%  "

# Row 45
$ws.Range("E45").Value = "  +5.82%  "

# Row 46
$ws.Range("E46").Value = "  +2.70%  "

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "96.51"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.48%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.87"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.01%  "

# Row 49
$ws.Range("D49").Value = "1.405.60"
$ws.Range("E49").Value = "  +2.37%  "

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.93"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.30%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "47.80"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.42%  "

